$wb = $excel.ActiveWorkbook

# Sheet 1 = 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 1101
$ws1.Range("F7").Value = 589
$ws1.Range("F8").Value = 1512
$ws1.Range("F10").Value = 1416
$ws1.Range("F11").Value = 3052
$ws1.Range("F12").Value = 557
$ws1.Range("F14").Value = 1782
$ws1.Range("F17").Value = 1445
$ws1.Range("F20").Value = 1177
$ws1.Range("F21").Value = 387
$ws1.Range("F22").Value = 429
$ws1.Range("F23").Value = 57
$ws1.Range("F24").Value = 4643
$ws1.Range("F25").Value = 729
$ws1.Range("F28").Value = 34
$ws1.Range("F29").Value = 82

# Sheet 2 = 演出 (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("G2").Value = "不可售"
$ws2.Range("F4").Value = 44
$ws2.Range("F7").Value = 19
$ws2.Range("F8").Value = 41
$ws2.Range("F10").Value = 23

# Sheet 4 = 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("G5").Value = "不可售"
$ws4.Range("F7").Value = 44
$ws4.Range("F11").Value = 19
$ws4.Range("F12").Value = 41
$ws4.Range("F14").Value = 23
$ws4.Range("F15").Value = 1101
$ws4.Range("F18").Value = 589
$ws4.Range("F19").Value = 1512
$ws4.Range("F21").Value = 1416
$ws4.Range("F22").Value = 3052
$ws4.Range("F23").Value = 557
$ws4.Range("F25").Value = 1782
$ws4.Range("F28").Value = 1445
$ws4.Range("F33").Value = 1177
$ws4.Range("F34").Value = 387
$ws4.Range("F35").Value = 429
$ws4.Range("F36").Value = 57
$ws4.Range("F37").Value = 4643
$ws4.Range("F38").Value = 729
$ws4.Range("F43").Value = 34
$ws4.Range("F44").Value = 82
